$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 9600
$ws.Range("F5").Value = 790
$ws.Range("F7").Value = 252
$ws.Range("G7").Value = 6.8
$ws.Range("F8").Value = 333
$ws.Range("F9").Value = 12
$ws.Range("F11").Value = 1446
$ws.Range("F15").Value = 319
$ws.Range("F19").Value = 419
$ws.Range("F25").Value = 293
$ws.Range("F27").Value = 274
$ws.Range("F35").Value = 57
$ws.Range("F36").Value = 192
$ws.Range("F37").Value = 346
$ws.Range("F39").Value = 347
$ws.Range("F46").Value = 64

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F23").Value = 668
$ws.Range("F24").Value = 41
$ws.Range("F25").Value = 8
$ws.Range("F32").Value = 129
$ws.Range("F33").Value = 176
$ws.Range("F34").Value = 33
$ws.Range("F35").Value = 26

# ---- Sheet "本地生活" ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 785
$ws.Range("F5").Value = 163
$ws.Range("F6").Value = 2434
$ws.Range("F7").Value = 3792
$ws.Range("F8").Value = 27
$ws.Range("F10").Value = 144
$ws.Range("F11").Value = 133

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 785
$ws.Range("F4").Value = 9600
$ws.Range("F5").Value = 163
$ws.Range("F7").Value = 3792
$ws.Range("F8").Value = 790
$ws.Range("F9").Value = 144
$ws.Range("F10").Value = 144
$ws.Range("F12").Value = 252
$ws.Range("G12").Value = 6.8
$ws.Range("F13").Value = 333
$ws.Range("F16").Value = 1446
$ws.Range("F18").Value = 133
$ws.Range("F21").Value = 319
$ws.Range("F28").Value = 293
$ws.Range("F31").Value = 274
$ws.Range("F35").Value = 41
$ws.Range("F36").Value = 8
$ws.Range("F38").Value = 57
$ws.Range("F39").Value = 346
$ws.Range("F41").Value = 347
$ws.Range("F48").Value = 332
$ws.Range("F49").Value = 331
